$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new JAN-2021 worksheet, placed right after DEC-2020 (last sheet)
# ---------------------------------------------------------------------------
$decSheet = $wb.Worksheets.Item("DEC-2020")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $decSheet)
$newSheet.Name = "JAN-2021"

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$newSheet.Cells.Item(1, 1).Value = "No"
$newSheet.Cells.Item(1, 2).Value = "Date"
$newSheet.Cells.Item(1, 3).Value = "Application"
$newSheet.Cells.Item(1, 4).Value = "Task"
$newSheet.Cells.Item(1, 5).Value = "% of completion"
$newSheet.Cells.Item(1, 6).Value = "Status"
$newSheet.Cells.Item(1, 7).Value = "Comments"

$headerRange = $newSheet.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0.39997558519241921
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------------
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = "1/1/2021"
$newSheet.Cells.Item(2, 4).Value = "Holiday"

$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = "1/2/2021"
$newSheet.Cells.Item(3, 4).Value = "Week off"

$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = "1/3/2021"
$newSheet.Cells.Item(4, 4).Value = "Week off"

$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = "1/4/2021"
$newSheet.Cells.Item(5, 3).Value = "B2B app & Qmvar 2.0"
$newSheet.Cells.Item(5, 4).Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. `nRegression testing on B2B app.`nRegression testing, Retesting and Cross browser testing on QMVAR 2.0 application"
$newSheet.Cells.Item(5, 5).Value = 1
$newSheet.Cells.Item(5, 6).Value = "Completed"
$newSheet.Rows.Item(5).RowHeight = 60

# Borders around whole data body
$bodyRange = $newSheet.Range("A1:G5")
$bodyRange.Borders.LineStyle = 1

# Date column formatting
$dateRange = $newSheet.Range("B2:B5")
$dateRange.NumberFormat = "[$-14009]yyyy-mm-dd;@"
$dateRange.HorizontalAlignment = -4131

# Task/Application columns -> wrap text, left aligned
$wrapRange = $newSheet.Range("C2:D5")
$wrapRange.WrapText = $true
$wrapRange.HorizontalAlignment = -4131

# Holiday / Week off (D2:D4) centered, bold red, wrap
$offRange = $newSheet.Range("D2:D4")
$offRange.HorizontalAlignment = -4108
$offRange.Font.Bold = $true
$offRange.Font.Color = 255

# % of completion column
$pctRange = $newSheet.Range("E2:E5")
$pctRange.NumberFormat = "0%"
$pctRange.HorizontalAlignment = -4131

# Status column - "Completed" cell gets green fill
$statusCell = $newSheet.Cells.Item(5, 6)
$statusCell.Interior.Pattern = 1
$statusCell.Interior.Color = 5287936
$statusCell.HorizontalAlignment = -4131

# Generic left alignment for remaining body cells
$newSheet.Range("A2:A5").HorizontalAlignment = -4131
$newSheet.Range("F2:F5").HorizontalAlignment = -4131
$newSheet.Range("G1:G5").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 4. Column widths (matching sibling monthly sheets)
# ---------------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 3.5703125
$newSheet.Columns.Item(2).ColumnWidth = 10.42578125
$newSheet.Columns.Item(3).ColumnWidth = 11.140625
$newSheet.Columns.Item(4).ColumnWidth = 74.28515625
$newSheet.Columns.Item(5).ColumnWidth = 15.5703125
$newSheet.Columns.Item(6).ColumnWidth = 10.85546875
$newSheet.Columns.Item(7).ColumnWidth = 10.5703125

# ---------------------------------------------------------------------------
# 5. View state: make JAN-2021 the active/selected tab
# ---------------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("D6").Select()

# DEC-2020 keeps a plain (non-active) view with an updated selection
$decSheet.Range("D32").Select()

$newSheet.Activate()
